# Apply cryptos list update (prices/volume refresh + Arweave/TheGraph row swap)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.285.21"
$ws.Range("E2").Value = "  +7.68%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.018.55"
$ws.Range("E3").Value = "  +4.81%  "

# Row 4
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "583.65"
$ws.Range("E5").Value = "  +3.19%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.85"
$ws.Range("E6").Value = "  +10.24%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.016.56"
$ws.Range("E8").Value = "  +4.84%  "

# Row 9
$ws.Range("E9").Value = "  +3.71%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.99"
$ws.Range("E10").Value = "  +1.86%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.156"
$ws.Range("E11").Value = "  +7.72%  "

# Row 12
$ws.Range("E12").Value = "  +5.84%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000254"
$ws.Range("E13").Value = "  +10.69%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.57"
$ws.Range("E14").Value = "  +9.69%  "

# Row 15
$ws.Range("E15").Value = "  +0.74%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.214.89"
$ws.Range("E16").Value = "  +7.63%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.517.41"
$ws.Range("E17").Value = "  +4.74%  "

# Row 18
$ws.Range("E18").Value = "  +7.32%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.021.84"
$ws.Range("E19").Value = "  +4.67%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "463.84"
$ws.Range("E20").Value = "  +8.33%  "

# Row 21
$ws.Range("E21").Value = "  +7.27%  "

# Row 22
$ws.Range("E22").Value = "  +5.45%  "

# Row 23
$ws.Range("E23").Value = "  +8.69%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "82.25"
$ws.Range("E24").Value = "  +4.47%  "

# Row 25
$ws.Range("E25").Value = "  +13.28%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.48"
$ws.Range("E26").Value = "  +5.70%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.68"
$ws.Range("E27").Value = "  +7.98%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.06"
$ws.Range("E29").Value = "  +14.75%  "

# Row 30
$ws.Range("E30").Value = "  +17.82%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0000106"
$ws.Range("E31").Value = "  +1.57%  "

# Row 32
$ws.Range("E32").Value = "  +5.19%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "27.06"
$ws.Range("E33").Value = "  +6.90%  "

# Row 34
$ws.Range("E34").Value = "  +5.33%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.999"
$ws.Range("E35").Value = "  -0.17%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.994"
$ws.Range("E36").Value = "  +4.37%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.79"
$ws.Range("E37").Value = "  +8.66%  "

# Row 38
$ws.Range("E38").Value = "  +14.47%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.07"
$ws.Range("E39").Value = "  +10.06%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "49.57"
$ws.Range("E40").Value = "  +1.64%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.123"
$ws.Range("E41").Value = "  +8.57%  "

# Row 42
$ws.Range("B42").Value = "Arweave"
$ws.Range("C42").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "43.83"
$ws.Range("E42").Value = "  +11.34%  "

# Row 43
$ws.Range("B43").Value = "TheGraph"
$ws.Range("C43").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.302"
$ws.Range("E43").Value = "  +14.47%  "

# Row 44
$ws.Range("E44").Value = "  +3.73%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "391.15"
$ws.Range("E45").Value = "  +14.54%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.806.14"
$ws.Range("E46").Value = "  +4.90%  "

# Row 47
$ws.Range("E47").Value = "  +6.39%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "133.96"
$ws.Range("E48").Value = "  +1.21%  "

# Row 50
$ws.Range("E50").Value = "  +10.41%  "

# Row 51
$ws.Range("E51").Value = "  +4.72%  "
